$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(298, 290, 183, 177, 169, 156, 136, 92, 90, 55, 42)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
